# Auto-generated edit script: refresh market-price-derived columns (H-N)
# on the per-job profit sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1570.7931
$ws.Range("I15").Value = 1570.7931
$ws.Range("K15").Value = 4712.379300000001
$ws.Range("M15").Value = -4543.379300000001

$ws.Range("H80").Value = 6090.8945
$ws.Range("I80").Value = 306.44446
$ws.Range("J80").Value = 11296.9
$ws.Range("K80").Value = 919.33338
$ws.Range("L80").Value = 33890.7
$ws.Range("M80").Value = 78.66661999999997
$ws.Range("N80").Value = -35886.7

$ws.Range("H83").Value = 6090.8945
$ws.Range("I83").Value = 306.44446
$ws.Range("J83").Value = 11296.9
$ws.Range("K83").Value = 2758.00014
$ws.Range("L83").Value = 101672.1
$ws.Range("M83").Value = 2233.99986
$ws.Range("N83").Value = -111656.1

$ws.Range("H98").Value = 5869.2583
$ws.Range("I98").Value = 3401.7083
$ws.Range("J98").Value = 14329.429
$ws.Range("K98").Value = 3401.7083
$ws.Range("L98").Value = 14329.429
$ws.Range("M98").Value = -1903.7083
$ws.Range("N98").Value = -17325.429

$ws.Range("H113").Value = 2701.111
$ws.Range("I113").Value = 2362
$ws.Range("J113").Value = 3125
$ws.Range("K113").Value = 2362
$ws.Range("L113").Value = 3125
$ws.Range("M113").Value = 892
$ws.Range("N113").Value = -9633

$ws.Range("H122").Value = 5869.2583
$ws.Range("I122").Value = 3401.7083
$ws.Range("J122").Value = 14329.429
$ws.Range("K122").Value = 10205.1249
$ws.Range("L122").Value = 42988.287
$ws.Range("M122").Value = -7755.124899999999
$ws.Range("N122").Value = -47888.287

$ws.Range("H127").Value = 1162.4
$ws.Range("I127").Value = 375.8
$ws.Range("K127").Value = 1127.4
$ws.Range("M127").Value = 3832.6

$ws.Range("H129").Value = 842.08
$ws.Range("I129").Value = 340.9375
$ws.Range("J129").Value = 937.5357
$ws.Range("K129").Value = 1022.8125
$ws.Range("L129").Value = 2812.6071
$ws.Range("M129").Value = 3977.1875
$ws.Range("N129").Value = -12812.6071

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1647
$ws.Range("I45").Value = 1092.8
$ws.Range("J45").Value = 2042.8572
$ws.Range("K45").Value = 1092.8
$ws.Range("L45").Value = 2042.8572
$ws.Range("M45").Value = -715.8
$ws.Range("N45").Value = -2796.8572

$ws.Range("H74").Value = 3353.1936
$ws.Range("I74").Value = 3357
$ws.Range("J74").Value = 3349.625
$ws.Range("K74").Value = 3357
$ws.Range("L74").Value = 3349.625
$ws.Range("M74").Value = -2483
$ws.Range("N74").Value = -5097.625

$ws.Range("H77").Value = 3353.1936
$ws.Range("I77").Value = 3357
$ws.Range("J77").Value = 3349.625
$ws.Range("K77").Value = 16785
$ws.Range("L77").Value = 16748.125
$ws.Range("M77").Value = -12417
$ws.Range("N77").Value = -25484.125

$ws.Range("H92").Value = 119058.336
$ws.Range("J92").Value = 119058.336
$ws.Range("L92").Value = 119058.336
$ws.Range("N92").Value = -124050.336

$ws.Range("H122").Value = 10604.826
$ws.Range("I122").Value = 10632.318
$ws.Range("K122").Value = 31896.954
$ws.Range("M122").Value = -29446.954

$ws.Range("H131").Value = 53142
$ws.Range("J131").Value = 53142
$ws.Range("L131").Value = 53142
$ws.Range("N131").Value = -63222

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1427.9333
$ws.Range("I99").Value = 1535.0834
$ws.Range("J99").Value = 999.3333
$ws.Range("K99").Value = 1535.0834
$ws.Range("L99").Value = 999.3333
$ws.Range("M99").Value = -37.08339999999998
$ws.Range("N99").Value = -3995.3333

$ws.Range("H107").Value = 24806.738
$ws.Range("I107").Value = 44348.832
$ws.Range("J107").Value = 3488.0908
$ws.Range("K107").Value = 44348.832
$ws.Range("L107").Value = 3488.0908
$ws.Range("M107").Value = -42428.832
$ws.Range("N107").Value = -7328.0908

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4692.8525
$ws.Range("I31").Value = 2520
$ws.Range("J31").Value = 5118.902
$ws.Range("K31").Value = 2520
$ws.Range("L31").Value = 5118.902
$ws.Range("M31").Value = -2225
$ws.Range("N31").Value = -5708.902

$ws.Range("H34").Value = 4692.8525
$ws.Range("I34").Value = 2520
$ws.Range("J34").Value = 5118.902
$ws.Range("K34").Value = 2520
$ws.Range("L34").Value = 5118.902
$ws.Range("M34").Value = -2318
$ws.Range("N34").Value = -5522.902

$ws.Range("H99").Value = 2283.0908
$ws.Range("I99").Value = 2328.5715
$ws.Range("J99").Value = 2203.5
$ws.Range("K99").Value = 2328.5715
$ws.Range("L99").Value = 2203.5
$ws.Range("M99").Value = -830.5715
$ws.Range("N99").Value = -5199.5

$ws.Range("H126").Value = 2283.0908
$ws.Range("I126").Value = 2328.5715
$ws.Range("J126").Value = 2203.5
$ws.Range("K126").Value = 6985.7145
$ws.Range("L126").Value = 6610.5
$ws.Range("M126").Value = -4515.7145
$ws.Range("N126").Value = -11550.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 8583
$ws.Range("I17").Value = 9000
$ws.Range("J17").Value = 8545.091
$ws.Range("K17").Value = 27000
$ws.Range("L17").Value = 25635.273
$ws.Range("M17").Value = -26831
$ws.Range("N17").Value = -25973.273

$ws.Range("H34").Value = 1248.5
$ws.Range("I34").Value = 188.88889
$ws.Range("J34").Value = 3155.8
$ws.Range("K34").Value = 566.6666700000001
$ws.Range("L34").Value = 9467.400000000001
$ws.Range("M34").Value = -482.6666700000001
$ws.Range("N34").Value = -9635.400000000001

$ws.Range("H39").Value = 4963.9165
$ws.Range("J39").Value = 4963.9165
$ws.Range("L39").Value = 14891.7495
$ws.Range("N39").Value = -15479.7495

$ws.Range("H55").Value = 6099.125
$ws.Range("I55").Value = 1000
$ws.Range("J55").Value = 6827.5713
$ws.Range("K55").Value = 3000
$ws.Range("L55").Value = 20482.7139
$ws.Range("M55").Value = -2823
$ws.Range("N55").Value = -20836.7139

$ws.Range("H98").Value = 101570.3
$ws.Range("I98").Value = 168617.17
$ws.Range("K98").Value = 505851.51
$ws.Range("M98").Value = -504353.51

$ws.Range("H109").Value = 10482.077
$ws.Range("I109").Value = 26481.75
$ws.Range("J109").Value = 3371.111
$ws.Range("K109").Value = 79445.25
$ws.Range("L109").Value = 10113.333
$ws.Range("M109").Value = -78405.25
$ws.Range("N109").Value = -12193.333

$ws.Range("H131").Value = 22213.535
$ws.Range("I131").Value = 608.4
$ws.Range("J131").Value = 28760.545
$ws.Range("K131").Value = 1825.2
$ws.Range("L131").Value = 86281.63499999999
$ws.Range("M131").Value = 3214.8
$ws.Range("N131").Value = -96361.63499999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("N5").ClearContents()

$ws.Range("H80").Value = 3000.6667
$ws.Range("J80").Value = 3503
$ws.Range("L80").Value = 3503
$ws.Range("N80").Value = -5499

$ws.Range("H83").Value = 3000.6667
$ws.Range("J83").Value = 3503
$ws.Range("L83").Value = 17515
$ws.Range("N83").Value = -27499

$ws.Range("H102").Value = 3124.8408
$ws.Range("I102").Value = 3033.658
$ws.Range("J102").Value = 3702.3333
$ws.Range("K102").Value = 3033.658
$ws.Range("L102").Value = 3702.3333
$ws.Range("M102").Value = -1411.658
$ws.Range("N102").Value = -6946.3333

$ws.Range("H122").Value = 3758.6667
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 3758.6667
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 11276.0001
$ws.Range("N122").Value = -16176.0001
$ws.Range("M122").ClearContents()

$ws.Range("H126").Value = 4071.0557
$ws.Range("I126").Value = 2054
$ws.Range("J126").Value = 5079.5835
$ws.Range("K126").Value = 6162
$ws.Range("L126").Value = 15238.7505
$ws.Range("M126").Value = -3692
$ws.Range("N126").Value = -20178.7505

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3720.1765
$ws.Range("I7").Value = 3116.4443
$ws.Range("J7").Value = 4399.375
$ws.Range("K7").Value = 3116.4443
$ws.Range("L7").Value = 4399.375
$ws.Range("M7").Value = -3004.4443
$ws.Range("N7").Value = -4623.375

$ws.Range("H22").Value = 1456.7778
$ws.Range("I22").Value = 2168
$ws.Range("J22").Value = 1183.2307
$ws.Range("K22").Value = 2168
$ws.Range("L22").Value = 1183.2307
$ws.Range("M22").Value = -1873
$ws.Range("N22").Value = -1773.2307

$ws.Range("H27").Value = 1456.7778
$ws.Range("I27").Value = 2168
$ws.Range("J27").Value = 1183.2307
$ws.Range("K27").Value = 2168
$ws.Range("L27").Value = 1183.2307
$ws.Range("M27").Value = -2061
$ws.Range("N27").Value = -1397.2307

$ws.Range("H40").Value = 3495.7334
$ws.Range("I40").Value = 3487.2
$ws.Range("K40").Value = 3487.2
$ws.Range("M40").Value = -3351.2

$ws.Range("H43").Value = 50000
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").ClearContents()

$ws.Range("H126").Value = 3720.1765
$ws.Range("I126").Value = 3116.4443
$ws.Range("J126").Value = 4399.375
$ws.Range("K126").Value = 9349.332900000001
$ws.Range("L126").Value = 13198.125
$ws.Range("M126").Value = -6879.332900000001
$ws.Range("N126").Value = -18138.125

$ws.Range("H132").Value = 3996.926
$ws.Range("I132").Value = 4750.5
$ws.Range("J132").Value = 2900.818
$ws.Range("K132").Value = 14251.5
$ws.Range("L132").Value = 8702.454000000002
$ws.Range("M132").Value = -11721.5
$ws.Range("N132").Value = -13762.454

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H123").Value = 23599.29
$ws.Range("J123").Value = 23599.29
$ws.Range("L123").Value = 23599.29
$ws.Range("N123").Value = -33399.29

$ws.Range("H125").Value = 49947
$ws.Range("J125").Value = 49947
$ws.Range("L125").Value = 49947
$ws.Range("N125").Value = -59787

$ws.Range("H126").Value = 4357.1177
$ws.Range("I126").Value = 6123.25
$ws.Range("K126").Value = 18369.75
$ws.Range("M126").Value = -15899.75

$ws.Range("H131").Value = 59975
$ws.Range("J131").Value = 59975
$ws.Range("L131").Value = 59975
$ws.Range("N131").Value = -70055

$ws.Range("H132").Value = 3937.2593
$ws.Range("I132").Value = 3526.8948
$ws.Range("J132").Value = 4911.875
$ws.Range("K132").Value = 10580.6844
$ws.Range("L132").Value = 14735.625
$ws.Range("M132").Value = -8050.6844
$ws.Range("N132").Value = -19795.625
